$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update progress for "LoiCamOn" (row 3): mark as finished / pages done
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "YES"
$ws.Range("E3").Value = "DONE"

# Narrow column D now that the long "ALMOST THERE..." text no longer needs the extra width
$ws.Columns.Item(4).ColumnWidth = 9.667

# Move the active selection
$ws.Range("C7").Select()

$wb.Save()
